# Update "all http result" values across the three result sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "Transfer Time (s)" ---
$ws = $wb.Worksheets.Item("Transfer Time (s)")
$ws.Range("B4").Value = 0.009963054299354554
$ws.Range("C4").Value = 0.005181819649484069
$ws.Range("D4").Value = 0.01418705487251282
$ws.Range("E4").Value = 0.01797109972563448

$ws.Range("B5").Value = 0.01323233366012573
$ws.Range("C5").Value = 0.005749628294880777
$ws.Range("D5").Value = 0.04627740144729614
$ws.Range("E5").Value = 0.01178264574478354

$ws.Range("B6").Value = 0.02930171489715576
$ws.Range("C6").Value = 0.01060076133917406
$ws.Range("D6").Value = 0.1010946869850159
$ws.Range("E6").Value = 0.02983875185514805

$ws.Range("B7").Value = 0.1984585523605347
$ws.Range("D7").Value = 2.159919381141663

# --- Sheet "Throughput (bps)" ---
$ws = $wb.Worksheets.Item("Throughput (bps)")
$ws.Range("B4").Value = 9827463.639666218
$ws.Range("C4").Value = 3817495.177306633
$ws.Range("D4").Value = 20321686.671219
$ws.Range("E4").Value = 14653429.31618624

$ws.Range("B5").Value = 73611763.66876601
$ws.Range("C5").Value = 30069905.29991582
$ws.Range("D5").Value = 22454243.5859525
$ws.Range("E5").Value = 21727084.20517864

$ws.Range("B6").Value = 324711003.9069602
$ws.Range("C6").Value = 114405569.2885726
$ws.Range("D6").Value = 86919085.30927305
$ws.Range("E6").Value = 15208536.58317422

$ws.Range("B7").Value = 425548546.5109209
$ws.Range("D7").Value = 38837619.61538628

# --- Sheet "Overhead Ratio" ---
$ws = $wb.Worksheets.Item("Overhead Ratio")
$ws.Range("B4").Value = 1.0396484375
$ws.Range("B5").Value = 1.003984375
$ws.Range("B6").Value = 1.000388145446777
$ws.Range("B7").Value = 1.000039005279541
